# Add the new "fourth state" data row (row 4) under the existing A1:D3 block.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 1

# Reflect the author's last on-screen view: zoomed to 105% with C7 selected.
$ws.Range("C7").Select()
$excel.ActiveWindow.Zoom = 105
